$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3459.4119
$ws.Range("I64").Value = 3060
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 3060
$ws.Range("L64").Value = 5000
$ws.Range("M64").Value = -2812
$ws.Range("N64").Value = -5496

$ws.Range("H67").Value = 3459.4119
$ws.Range("I67").Value = 3060
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 3060
$ws.Range("L67").Value = 5000
$ws.Range("M67").Value = -2202
$ws.Range("N67").Value = -6716

$ws.Range("H113").Value = 8594.522999999999
$ws.Range("I113").Value = 2700
$ws.Range("J113").Value = 11541.786
$ws.Range("K113").Value = 2700
$ws.Range("L113").Value = 11541.786
$ws.Range("M113").Value = 554
$ws.Range("N113").Value = -18049.786

$ws.Range("H132").Value = 15183.933
$ws.Range("I132").Value = 16299.044
$ws.Range("J132").Value = 2546
$ws.Range("K132").Value = 48897.132
$ws.Range("L132").Value = 7638
$ws.Range("M132").Value = -46367.132
$ws.Range("N132").Value = -12698

$ws.Range("H138").Value = 6157.91
$ws.Range("I138").Value = 4149.077
$ws.Range("J138").Value = 6559.677
$ws.Range("K138").Value = 12447.231
$ws.Range("L138").Value = 19679.031
$ws.Range("M138").Value = -7307.231
$ws.Range("N138").Value = -29959.031

$ws.Range("H141").Value = 2834.5
$ws.Range("I141").Value = 1618.9375
$ws.Range("K141").Value = 4856.8125
$ws.Range("M141").Value = 323.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15116.742
$ws.Range("I32").Value = 11592.087
$ws.Range("J32").Value = 25250.125
$ws.Range("K32").Value = 11592.087
$ws.Range("L32").Value = 25250.125
$ws.Range("M32").Value = -11305.087
$ws.Range("N32").Value = -25824.125

$ws.Range("H61").Value = 1989.1316
$ws.Range("I61").Value = 1736.037
$ws.Range("J61").Value = 2610.3635
$ws.Range("K61").Value = 1736.037
$ws.Range("L61").Value = 2610.3635
$ws.Range("M61").Value = -1524.037
$ws.Range("N61").Value = -3034.3635

$ws.Range("H114").Value = 37975
$ws.Range("J114").Value = 37975
$ws.Range("L114").Value = 37975
$ws.Range("N114").Value = -46653

$ws.Range("H132").Value = 1522.9445
$ws.Range("I132").Value = 1145.3091
$ws.Range("J132").Value = 2744.7058
$ws.Range("K132").Value = 3435.9273
$ws.Range("L132").Value = 8234.117400000001
$ws.Range("M132").Value = -905.9272999999998
$ws.Range("N132").Value = -13294.1174

$ws.Range("H136").Value = 1989.1316
$ws.Range("I136").Value = 1736.037
$ws.Range("J136").Value = 2610.3635
$ws.Range("K136").Value = 5208.111
$ws.Range("L136").Value = 7831.0905
$ws.Range("M136").Value = -2658.111
$ws.Range("N136").Value = -12931.0905

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7435.4
$ws.Range("I20").Value = 11671.091
$ws.Range("J20").Value = 2258.4443
$ws.Range("K20").Value = 11671.091
$ws.Range("L20").Value = 2258.4443
$ws.Range("M20").Value = -11424.091
$ws.Range("N20").Value = -2752.4443

$ws.Range("H94").Value = 1012.5
$ws.Range("I94").Value = 1012.5
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1012.5
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -561.5
$ws.Range("N94").ClearContents()

$ws.Range("H132").Value = 46783.332
$ws.Range("J132").Value = 46783.332
$ws.Range("L132").Value = 46783.332
$ws.Range("N132").Value = -56903.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2825
$ws.Range("J22").Value = 433.33334
$ws.Range("L22").Value = 433.33334
$ws.Range("N22").Value = -1133.33334

$ws.Range("H132").Value = 2226.1853
$ws.Range("I132").Value = 1175.25
$ws.Range("J132").Value = 3754.818
$ws.Range("K132").Value = 3525.75
$ws.Range("L132").Value = 11264.454
$ws.Range("M132").Value = -995.75
$ws.Range("N132").Value = -16324.454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 3682.2
$ws.Range("J64").Value = 4999.6665
$ws.Range("L64").Value = 14998.9995
$ws.Range("N64").Value = -15538.9995

$ws.Range("H67").Value = 3682.2
$ws.Range("J67").Value = 4999.6665
$ws.Range("L67").Value = 14998.9995
$ws.Range("N67").Value = -16870.9995

$ws.Range("H113").Value = 1070.4546
$ws.Range("I113").Value = 1214.4375
$ws.Range("J113").Value = 686.5
$ws.Range("K113").Value = 3643.3125
$ws.Range("L113").Value = 2059.5
$ws.Range("M113").Value = -1473.3125
$ws.Range("N113").Value = -6399.5

$ws.Range("H122").Value = 1079.1852
$ws.Range("I122").Value = 619.087
$ws.Range("K122").Value = 5571.782999999999
$ws.Range("M122").Value = -3121.782999999999

$ws.Range("H137").Value = 8344.807000000001
$ws.Range("I137").Value = 2509.889
$ws.Range("J137").Value = 10731.818
$ws.Range("K137").Value = 7529.667
$ws.Range("L137").Value = 32195.454
$ws.Range("M137").Value = -2429.667
$ws.Range("N137").Value = -42395.454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2161.6
$ws.Range("I126").Value = 1975.6923
$ws.Range("J126").Value = 2363
$ws.Range("K126").Value = 5927.0769
$ws.Range("L126").Value = 7089
$ws.Range("M126").Value = -3457.0769
$ws.Range("N126").Value = -12029

$ws.Range("H132").Value = 2594.7073
$ws.Range("I132").Value = 2252.3667
$ws.Range("J132").Value = 3528.3635
$ws.Range("K132").Value = 6757.1001
$ws.Range("L132").Value = 10585.0905
$ws.Range("M132").Value = -4227.1001
$ws.Range("N132").Value = -15645.0905

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1476.4615
$ws.Range("I22").Value = 999.1429000000001
$ws.Range("J22").Value = 2033.3334
$ws.Range("K22").Value = 999.1429000000001
$ws.Range("L22").Value = 2033.3334
$ws.Range("M22").Value = -704.1429000000001
$ws.Range("N22").Value = -2623.3334

$ws.Range("H27").Value = 1476.4615
$ws.Range("I27").Value = 999.1429000000001
$ws.Range("J27").Value = 2033.3334
$ws.Range("K27").Value = 999.1429000000001
$ws.Range("L27").Value = 2033.3334
$ws.Range("M27").Value = -892.1429000000001
$ws.Range("N27").Value = -2247.3334

$ws.Range("H68").Value = 2172.625
$ws.Range("J68").Value = 2172.625
$ws.Range("L68").Value = 2172.625
$ws.Range("N68").Value = -3670.625

$ws.Range("H71").Value = 2172.625
$ws.Range("J71").Value = 2172.625
$ws.Range("L71").Value = 10863.125
$ws.Range("N71").Value = -18351.125

$ws.Range("H122").Value = 10105053
$ws.Range("I122").Value = 15877272
$ws.Range("J122").Value = 3670
$ws.Range("K122").Value = 47631816
$ws.Range("L122").Value = 11010
$ws.Range("M122").Value = -47629366
$ws.Range("N122").Value = -15910

$ws.Range("H137").Value = 34132.715
$ws.Range("I137").Value = 13000
$ws.Range("J137").Value = 42585.8
$ws.Range("K137").Value = 13000
$ws.Range("L137").Value = 42585.8
$ws.Range("M137").Value = -7900
$ws.Range("N137").Value = -52785.8

$ws.Range("H139").Value = 44625
$ws.Range("J139").Value = 44625
$ws.Range("L139").Value = 44625
$ws.Range("N139").Value = -54905

$ws.Range("H141").Value = 67858.92999999999
$ws.Range("J141").Value = 67858.92999999999
$ws.Range("L141").Value = 67858.92999999999
$ws.Range("N141").Value = -78218.92999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 35349.965
$ws.Range("I122").Value = 45147.78
$ws.Range("J122").Value = 3157.1428
$ws.Range("K122").Value = 135443.34
$ws.Range("L122").Value = 9471.428400000001
$ws.Range("M122").Value = -132993.34
$ws.Range("N122").Value = -14371.4284

$ws.Range("H132").Value = 1540
$ws.Range("I132").Value = 1085.0339
$ws.Range("J132").Value = 2658.4583
$ws.Range("K132").Value = 3255.1017
$ws.Range("L132").Value = 7975.374899999999
$ws.Range("M132").Value = -725.1016999999997
$ws.Range("N132").Value = -13035.3749

$ws.Range("H136").Value = 5092.1333
$ws.Range("I136").Value = 3409.2415
$ws.Range("J136").Value = 6666.4517
$ws.Range("K136").Value = 10227.7245
$ws.Range("L136").Value = 19999.3551
$ws.Range("M136").Value = -7677.7245
$ws.Range("N136").Value = -25099.3551
